$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at the very left (A and B), shifting the existing
# A:G columns to C:I.
$ws.Range("A:B").EntireColumn.Insert()

# Copy the header formatting from the (now shifted) header cell C1 onto the
# two freshly inserted header cells so they match the bold/centered/bordered
# style used by the rest of row 1.
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

# New header labels for the inserted columns.
$ws.Range("A1").Value = "data"
$ws.Range("B1").Value = "loja"

# New tracking id used in every "link" URL (column I).
$oldTrackingId = "f58cbb8f-93eb-4efd-b58f-a6d1e4dca198"
$newTrackingId = "0665459a-e415-4da0-a7c8-7b68af51a563"

# Fill in the new "data" / "loja" values for every data row, and refresh the
# tracking_id query parameter embedded in the link column (now column I).
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = "30/07/2024"
    $ws.Cells.Item($r, 2).Value = "acessorios web"

    $link = $ws.Cells.Item($r, 9).Value2()
    $newLink = $link.Replace($oldTrackingId, $newTrackingId)
    $ws.Cells.Item($r, 9).Value = $newLink
}
